# Applies the cryptocurrency price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
# Column D ("Price") holds plain decimal-looking numbers that must stay as literal
# TEXT (matching the original inlineStr cells), so those are written with a leading
# apostrophe (backtick-escaped single quote) to force Excel to keep them as text
# instead of auto-converting to a numeric value (which would silently drop things
# like trailing zeros, e.g. "21.60" -> 21.6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.247.30"
$ws.Range("D3").Value = "1.906.67"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "`'307.65"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "`'0.5248"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "`'0.3813"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").Value = "`'0.07297"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "`'21.60"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "`'0.9049"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  -4.66%  "
$ws.Range("D13").Value = "`'96.42"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "`'5.366"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "1.599.80"
$ws.Range("E15").Value = "  -16.05%  "
$ws.Range("D17").Value = "`'0.000008677"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "`'14.76"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "27.069.63"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "`'5.124"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").Value = "`'10.81"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").Value = "`'6.514"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "`'2.357"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").Value = "`'149.83"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").Value = "`'18.26"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "`'1.735"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "`'116.67"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").Value = "`'4.853"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "`'4.866"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "`'0.09248"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").Value = "`'0.8355"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("D33").Value = "`'0.05074"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "`'1.231"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "`'2.983"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "`'3.356"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("D37").Value = "`'2.734"
$ws.Range("E37").Value = "  +4.52%  "
$ws.Range("D38").Value = "`'0.5770"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "`'0.02006"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "`'1.080"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("D41").Value = "`'9.141"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "`'6.584"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "`'116.42"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "`'0.1522"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "`'0.4922"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").Value = "`'10.21"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "`'1.001"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "`'1.644"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").Value = "`'38.78"
$ws.Range("E49").Value = "  +3.11%  "
$ws.Range("D50").Value = "`'64.51"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "`'0.06055"
$ws.Range("E51").Value = "  +1.72%  "
